$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2238442822384428
$ws.Range("C2").Value = 0.4987834549878346
$ws.Range("J2").Value = 0.0218978102189781
$ws.Range("O2").Value = 0.004866180048661801
$ws.Range("P2").Value = 0.145985401459854
$ws.Range("S2").Value = 0.1046228710462287
$ws.Range("B3").Value = 0.004651162790697674
$ws.Range("C3").Value = 0.02325581395348837
$ws.Range("J3").Value = 0.04186046511627907
$ws.Range("P3").Value = 0.7581395348837209
$ws.Range("S3").Value = 0.172093023255814
$ws.Range("J4").Value = 0.04347826086956522
$ws.Range("P4").Value = 0.5869565217391305
$ws.Range("S4").Value = 0.3695652173913043
$ws.Range("B6").Value = 0.09473684210526316
$ws.Range("D6").Value = 0.01578947368421053
$ws.Range("F6").Value = 0.01578947368421053
$ws.Range("J6").Value = 0.2789473684210526
$ws.Range("O6").Value = 0.01052631578947368
$ws.Range("Q6").Value = 0.1263157894736842
$ws.Range("R6").Value = 0.08947368421052632
$ws.Range("S6").Value = 0.3684210526315789
$ws.Range("B7").Value = 0.108695652173913
$ws.Range("D7").Value = 0.03260869565217391
$ws.Range("F7").Value = 0.05434782608695652
$ws.Range("J7").Value = 0.09782608695652174
$ws.Range("O7").Value = 0.01630434782608696
$ws.Range("Q7").Value = 0.2391304347826087
$ws.Range("R7").Value = 0.05978260869565218
$ws.Range("S7").Value = 0.391304347826087
$ws.Range("B8").Value = 0.1161290322580645
$ws.Range("D8").Value = 0.01935483870967742
$ws.Range("F8").Value = 0.04301075268817205
$ws.Range("J8").Value = 0.1204301075268817
$ws.Range("O8").Value = 0.01935483870967742
$ws.Range("Q8").Value = 0.2043010752688172
$ws.Range("R8").Value = 0.08602150537634409
$ws.Range("S8").Value = 0.3913978494623656
$ws.Range("B9").Value = 0.1066666666666667
$ws.Range("D9").Value = 0.006666666666666667
$ws.Range("F9").Value = 0.06666666666666667
$ws.Range("J9").Value = 0.1066666666666667
$ws.Range("O9").Value = 0.01333333333333333
$ws.Range("Q9").Value = 0.1866666666666667
$ws.Range("R9").Value = 0.1066666666666667
$ws.Range("S9").Value = 0.4066666666666667
$ws.Range("B10").Value = 0.1446453407510431
$ws.Range("D10").Value = 0.02016689847009736
$ws.Range("E10").Value = 0.0006954102920723226
$ws.Range("F10").Value = 0.06397774687065369
$ws.Range("J10").Value = 0.10778859527121
$ws.Range("O10").Value = 0.009040333796940195
$ws.Range("Q10").Value = 0.2329624478442281
$ws.Range("R10").Value = 0.05910987482614743
$ws.Range("S10").Value = 0.3616133518776078
$ws.Range("G11").Value = 0.1666666666666667
$ws.Range("J11").Value = 0.09477124183006536
$ws.Range("K11").Value = 0.2222222222222222
$ws.Range("L11").Value = 0.5
$ws.Range("S11").Value = 0.01633986928104575
$ws.Range("G12").Value = 0.6666666666666666
$ws.Range("J12").Value = 0.2592592592592592
$ws.Range("K12").Value = 0.01234567901234568
$ws.Range("L12").Value = 0.0308641975308642
$ws.Range("S12").Value = 0.0308641975308642
$ws.Range("G13").Value = 0.6428571428571429
$ws.Range("J13").Value = 0.3095238095238095
$ws.Range("S13").Value = 0.04761904761904762
$ws.Range("F15").Value = 0.009569377990430622
$ws.Range("H15").Value = 0.1913875598086124
$ws.Range("I15").Value = 0.05263157894736842
$ws.Range("J15").Value = 0.4449760765550239
$ws.Range("K15").Value = 0.06220095693779904
$ws.Range("M15").Value = 0.01435406698564593
$ws.Range("O15").Value = 0.06698564593301436
$ws.Range("S15").Value = 0.1578947368421053
$ws.Range("F16").Value = 0.00823045267489712
$ws.Range("H16").Value = 0.1522633744855967
$ws.Range("I16").Value = 0.07407407407407407
$ws.Range("J16").Value = 0.5102880658436214
$ws.Range("K16").Value = 0.07818930041152264
$ws.Range("M16").Value = 0.02469135802469136
$ws.Range("O16").Value = 0.01646090534979424
$ws.Range("S16").Value = 0.1358024691358025
$ws.Range("F17").Value = 0.01351351351351351
$ws.Range("H17").Value = 0.2007722007722008
$ws.Range("I17").Value = 0.06756756756756757
$ws.Range("J17").Value = 0.4555984555984556
$ws.Range("K17").Value = 0.0945945945945946
$ws.Range("M17").Value = 0.01158301158301158
$ws.Range("O17").Value = 0.05598455598455598
$ws.Range("S17").Value = 0.1003861003861004
$ws.Range("F18").Value = 0.01754385964912281
$ws.Range("H18").Value = 0.1695906432748538
$ws.Range("I18").Value = 0.05847953216374269
$ws.Range("J18").Value = 0.4970760233918128
$ws.Range("K18").Value = 0.08771929824561403
$ws.Range("M18").Value = 0.02339181286549707
$ws.Range("O18").Value = 0.07017543859649122
$ws.Range("S18").Value = 0.07602339181286549
$ws.Range("F19").Value = 0.01431980906921241
$ws.Range("H19").Value = 0.2084327764518695
$ws.Range("I19").Value = 0.06046141607000795
$ws.Range("J19").Value = 0.4073190135242641
$ws.Range("K19").Value = 0.1089896579156722
$ws.Range("M19").Value = 0.01988862370723946
$ws.Range("O19").Value = 0.07637231503579953
$ws.Range("S19").Value = 0.1042163882259348
